$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from E1 to F1, then set header value
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Fill time_taken values for data rows (no special style, matches columns B-E)
$ws.Range("F2").Value = "2021-10-05 10:52:07.630748"
$ws.Range("F3").Value = "2021-10-05 10:52:07.630760"
$ws.Range("F4").Value = "2021-10-05 10:52:07.630763"
$ws.Range("F5").Value = "2021-10-05 10:52:07.630766"
$ws.Range("F6").Value = "2021-10-05 10:52:07.630769"
$ws.Range("F7").Value = "2021-10-05 10:52:07.630771"
$ws.Range("F8").Value = "2021-10-05 10:52:07.630774"
$ws.Range("F9").Value = "2021-10-05 10:52:07.630777"
$ws.Range("F10").Value = "2021-10-05 10:52:07.630779"
$ws.Range("F11").Value = "2021-10-05 10:52:07.630782"
$ws.Range("F12").Value = "2021-10-05 10:52:07.630784"
$ws.Range("F13").Value = "2021-10-05 10:52:07.630787"
$ws.Range("F14").Value = "2021-10-05 10:52:07.630789"
$ws.Range("F15").Value = "2021-10-05 10:52:07.630792"
$ws.Range("F16").Value = "2021-10-05 10:52:07.630794"
$ws.Range("F17").Value = "2021-10-05 10:52:07.630797"
$ws.Range("F18").Value = "2021-10-05 10:52:07.630800"
$ws.Range("F19").Value = "2021-10-05 10:52:07.630802"
$ws.Range("F20").Value = "2021-10-05 10:52:07.630805"
$ws.Range("F21").Value = "2021-10-05 10:52:07.630807"
$ws.Range("F22").Value = "2021-10-05 10:52:07.630810"
$ws.Range("F23").Value = "2021-10-05 10:52:07.630812"
$ws.Range("F24").Value = "2021-10-05 10:52:07.630815"
$ws.Range("F25").Value = "2021-10-05 10:52:07.630817"
$ws.Range("F26").Value = "2021-10-05 10:52:07.630820"
$ws.Range("F27").Value = "2021-10-05 10:52:07.630823"
$ws.Range("F28").Value = "2021-10-05 10:52:07.630825"
$ws.Range("F29").Value = "2021-10-05 10:52:07.630828"
$ws.Range("F30").Value = "2021-10-05 10:52:07.630830"
$ws.Range("F31").Value = "2021-10-05 10:52:07.630833"
$ws.Range("F32").Value = "2021-10-05 10:52:07.630835"
$ws.Range("F33").Value = "2021-10-05 10:52:07.630838"
$ws.Range("F34").Value = "2021-10-05 10:52:07.630840"
$ws.Range("F35").Value = "2021-10-05 10:52:07.630843"
$ws.Range("F36").Value = "2021-10-05 10:52:07.630846"
$ws.Range("F37").Value = "2021-10-05 10:52:07.630848"
$ws.Range("F38").Value = "2021-10-05 10:52:07.630851"
$ws.Range("F39").Value = "2021-10-05 10:52:07.630853"
$ws.Range("F40").Value = "2021-10-05 10:52:07.630856"
$ws.Range("F41").Value = "2021-10-05 10:52:07.630858"
$ws.Range("F42").Value = "2021-10-05 10:52:07.630862"
$ws.Range("F43").Value = "2021-10-05 10:52:07.630864"
$ws.Range("F44").Value = "2021-10-05 10:52:07.630867"
$ws.Range("F45").Value = "2021-10-05 10:52:07.630869"
$ws.Range("F46").Value = "2021-10-05 10:52:07.630872"
$ws.Range("F47").Value = "2021-10-05 10:52:07.630875"
$ws.Range("F48").Value = "2021-10-05 10:52:07.630879"
$ws.Range("F49").Value = "2021-10-05 10:52:07.630881"
$ws.Range("F50").Value = "2021-10-05 10:52:07.630884"
$ws.Range("F51").Value = "2021-10-05 10:52:07.630886"
$ws.Range("F52").Value = "2021-10-05 10:52:07.630889"
$ws.Range("F53").Value = "2021-10-05 10:52:07.630891"
$ws.Range("F54").Value = "2021-10-05 10:52:07.630894"
$ws.Range("F55").Value = "2021-10-05 10:52:07.630897"
$ws.Range("F56").Value = "2021-10-05 10:52:07.630899"
$ws.Range("F57").Value = "2021-10-05 10:52:07.630902"
$ws.Range("F58").Value = "2021-10-05 10:52:07.630904"
$ws.Range("F59").Value = "2021-10-05 10:52:07.630906"
$ws.Range("F60").Value = "2021-10-05 10:52:07.630909"
$ws.Range("F61").Value = "2021-10-05 10:52:07.630911"
$ws.Range("F62").Value = "2021-10-05 10:52:07.630914"
$ws.Range("F63").Value = "2021-10-05 10:52:07.630916"
$ws.Range("F64").Value = "2021-10-05 10:52:07.630919"
$ws.Range("F65").Value = "2021-10-05 10:52:07.630922"

Write-Output "done"
